# The edit moves Word's special "_GoBack" bookmark (the "last edit
# location" marker) from where it used to sit - right after the text
# "Protocole P1 ; Médicaments A2 ANSM + Q1 CPP" - to the top of the
# document, right before the text "Protocole P1 + résumé ...".
#
# Adding a new "_GoBack" bookmark automatically removes/replaces the
# previous one (Word allows only a single "_GoBack" bookmark at a time),
# and every other bookmark's w:id is renumbered in document order, which
# is exactly the id shuffle described by the diff.

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("Protocole P1 + r")) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph 'Protocole P1 + résumé ...'"
}

$collapsed = $d.Range($target.Range.Start, $target.Range.Start)
$d.Bookmarks.Add("_GoBack", $collapsed) | Out-Null
